$d = $word.ActiveDocument

function Replace-Exact($findText, $replaceText) {
    $range = $d.Content
    $ok = $range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $findText"
    }
}

# 1. Siemens Technology and services pvt. ltd, Pune
Replace-Exact "Siemens Technology and services pvt. ltd, Pune" "Siemens Technology and services pvt. ltd, Pune"

# 2. Tata Elxsi ltd, Pune
Replace-Exact " Tata Elxsi ltd, Pune" " Tata Elxsi ltd, Pune"

# 3. Rockwell Collins India, pvt. Ltd, Hyderabad
Replace-Exact " Rockwell Collins India, pvt. Ltd, Hyderabad" " Rockwell Collins India, pvt. Ltd, Hyderabad"

# 4 & 5. Qcc, g++, gcc (appears twice, identical in both tables)
Replace-Exact "Qcc, g++, gcc" "Qcc, g++, gcc"
Replace-Exact "Qcc, g++, gcc" "Qcc, g++, gcc"

# 6. AppLink (standalone cell)
Replace-Exact "AppLink" "AppLink"

# 7. Applink inside the project description paragraph
Replace-Exact "then Applink involves" "then Applink involves"

# 8. Visual studio, SVN, Super tool, gcc  (merge) + append ", Makefile" after qt4.8
Replace-Exact "Visual studio, SVN, Super tool, gcc" "Visual studio, SVN, Super tool, gcc"

# 9. Append ", Makefile" after "qt4.8"
Replace-Exact "qt4.8" "qt4.8, Makefile"

# 10. Append ", Makefile" after "Visual studio, SVN, Doors, Clear quest, PREP"
Replace-Exact "Visual studio, SVN, Doors, Clear quest, PREP" "Visual studio, SVN, Doors, Clear quest, PREP, Makefile"

Write-Output "done"
